$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.37870000000003
$ws.Range("C5").Value = -14.0224
$ws.Range("D7").Value = -7.541699999999993
$ws.Range("C9").Value = -12.0096
$ws.Range("C11").Value = -13.1695
$ws.Range("D11").Value = -8.2849
$ws.Range("A21").Value = -21.08340000000001
$ws.Range("C21").Value = -10.6145
$ws.Range("D21").Value = -7.502299999999998
$ws.Range("A23").Value = -21.37530000000002
$ws.Range("A25").Value = -22.56600000000003
